# Added unit checks for check_dqoaccuracy #5
# Normalize the "uom" column to lowercase/simplified unit strings and
# rename the NO3 / NH3 parameter labels to their full names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accuracy")

# --- Unit of measure updates (column B), first occurrence order ---
$ws.Range("B4").Value  = "mg/l"
$ws.Range("B6").Value  = "uS/cm"
$ws.Range("B3").Value  = "s.u."
$ws.Range("B19").Value = "cfu/100ml"
$ws.Range("B17").Value = "ug/l"

# --- Parameter name updates (column A) ---
$ws.Range("A14").Value = "Nitrate"
$ws.Range("A15").Value = "Ammonia"

# --- Remaining unit of measure updates (column B) ---
$ws.Range("B5").Value  = "mg/l"
$ws.Range("B7").Value  = "uS/cm"
$ws.Range("B8").Value  = "mg/l"
$ws.Range("B9").Value  = "mg/l"
$ws.Range("B10").Value = "mg/l"
$ws.Range("B11").Value = "mg/l"
$ws.Range("B12").Value = "mg/l"
$ws.Range("B13").Value = "mg/l"
$ws.Range("B14").Value = "mg/l"
$ws.Range("B15").Value = "mg/l"
$ws.Range("B16").Value = "mg/l"
$ws.Range("B18").Value = "ug/l"
$ws.Range("B20").Value = "cfu/100ml"
$ws.Range("B21").Value = "cfu/100ml"
$ws.Range("B22").Value = "cfu/100ml"

# --- Leave the cursor where the author left it when saving ---
$ws.Range("A15").Select() | Out-Null
